# Auto-applies the market-data refresh described by the commit diff.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for
# specific rows across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 179.54546
$ws.Range("I2").Value = 179
$ws.Range("J2").Value = 185
$ws.Range("K2").Value = 179
$ws.Range("L2").Value = 185
$ws.Range("M2").Value = -66
$ws.Range("N2").Value = -411

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4153.8887
$ws.Range("I86").Value = 2997
$ws.Range("J86").Value = 5600
$ws.Range("K86").Value = 2997
$ws.Range("L86").Value = 5600
$ws.Range("M86").Value = -1874
$ws.Range("N86").Value = -7846

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2956.2856
$ws.Range("I88").Value = 3559.75
$ws.Range("J88").Value = 2151.6667
$ws.Range("K88").Value = 3559.75
$ws.Range("L88").Value = 2151.6667
$ws.Range("M88").Value = -3153.75
$ws.Range("N88").Value = -2963.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4153.8887
$ws.Range("I89").Value = 2997
$ws.Range("J89").Value = 5600
$ws.Range("K89").Value = 14985
$ws.Range("L89").Value = 28000
$ws.Range("M89").Value = -9369
$ws.Range("N89").Value = -39232

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2956.2856
$ws.Range("I91").Value = 3559.75
$ws.Range("J91").Value = 2151.6667
$ws.Range("K91").Value = 3559.75
$ws.Range("L91").Value = 2151.6667
$ws.Range("M91").Value = -2155.75
$ws.Range("N91").Value = -4959.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 9545
$ws.Range("J113").Value = 11800
$ws.Range("L113").Value = 11800
$ws.Range("N113").Value = -18308

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 9920.375
$ws.Range("I131").Value = 6892.933
$ws.Range("K131").Value = 20678.799
$ws.Range("M131").Value = -15638.799

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 31251968
$ws.Range("I135").Value = 33335232
$ws.Range("K135").Value = 300017088
$ws.Range("M135").Value = -300014553

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 79285800
$ws.Range("J5").Value = 200000020
$ws.Range("L5").Value = 200000020
$ws.Range("N5").Value = -200000244

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1560.26
$ws.Range("I32").Value = 1492.8667
$ws.Range("J32").Value = 2166.8
$ws.Range("K32").Value = 1492.8667
$ws.Range("L32").Value = 2166.8
$ws.Range("M32").Value = -1205.8667
$ws.Range("N32").Value = -2740.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4783.357
$ws.Range("I45").Value = 3107.5557
$ws.Range("J45").Value = 7799.8
$ws.Range("K45").Value = 3107.5557
$ws.Range("L45").Value = 7799.8
$ws.Range("M45").Value = -2730.5557
$ws.Range("N45").Value = -8553.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1714.973
$ws.Range("I132").Value = 1551.9354
$ws.Range("K132").Value = 4655.8062
$ws.Range("M132").Value = -2125.8062

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 79285800
$ws.Range("J4").Value = 200000020
$ws.Range("L4").Value = 200000020
$ws.Range("N4").Value = -200000250

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2672.1667
$ws.Range("I94").Value = 2261.1177
$ws.Range("J94").Value = 3670.4285
$ws.Range("K94").Value = 2261.1177
$ws.Range("L94").Value = 3670.4285
$ws.Range("M94").Value = -1810.1177
$ws.Range("N94").Value = -4572.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1133.3334
$ws.Range("J14").Value = 1220
$ws.Range("L14").Value = 1220
$ws.Range("N14").Value = -1560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1654.0714
$ws.Range("I31").Value = 1666
$ws.Range("J31").Value = 1624.25
$ws.Range("K31").Value = 1666
$ws.Range("L31").Value = 1624.25
$ws.Range("M31").Value = -1371
$ws.Range("N31").Value = -2214.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1654.0714
$ws.Range("I34").Value = 1666
$ws.Range("J34").Value = 1624.25
$ws.Range("K34").Value = 1666
$ws.Range("L34").Value = 1624.25
$ws.Range("M34").Value = -1464
$ws.Range("N34").Value = -2028.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 30119.625
$ws.Range("J86").Value = 19774.715
$ws.Range("L86").Value = 19774.715
$ws.Range("N86").Value = -22020.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 30119.625
$ws.Range("J89").Value = 19774.715
$ws.Range("L89").Value = 98873.575
$ws.Range("N89").Value = -110105.575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2910.7778
$ws.Range("J99").Value = 2401.875
$ws.Range("L99").Value = 2401.875
$ws.Range("N99").Value = -5397.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3811.0286
$ws.Range("I122").Value = 3169.6
$ws.Range("J122").Value = 4666.2666
$ws.Range("K122").Value = 9508.799999999999
$ws.Range("L122").Value = 13998.7998
$ws.Range("M122").Value = -7058.799999999999
$ws.Range("N122").Value = -18898.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2910.7778
$ws.Range("J126").Value = 2401.875
$ws.Range("L126").Value = 7205.625
$ws.Range("N126").Value = -12145.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2554.923
$ws.Range("J132").Value = 4998.75
$ws.Range("L132").Value = 14996.25
$ws.Range("N132").Value = -20056.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3349.0715
$ws.Range("I134").Value = 3322.8333
$ws.Range("K134").Value = 9968.499899999999
$ws.Range("M134").Value = -7433.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2008.7
$ws.Range("J113").Value = 1665.3334
$ws.Range("L113").Value = 4996.0002
$ws.Range("N113").Value = -9336.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 3790.8333
$ws.Range("I114").Value = 1100
$ws.Range("J114").Value = 5136.25
$ws.Range("K114").Value = 3300
$ws.Range("L114").Value = 15408.75
$ws.Range("M114").Value = -46
$ws.Range("N114").Value = -21916.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 5534.778
$ws.Range("I36").Value = 955.4
$ws.Range("J36").Value = 11259
$ws.Range("K36").Value = 955.4
$ws.Range("L36").Value = 11259
$ws.Range("M36").Value = -470.4
$ws.Range("N36").Value = -12229

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 8199.23
$ws.Range("I97").Value = 533.125
$ws.Range("J97").Value = 20465
$ws.Range("K97").Value = 533.125
$ws.Range("L97").Value = 20465
$ws.Range("M97").Value = -37.125
$ws.Range("N97").Value = -21457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2021.3
$ws.Range("I126").Value = 1928.4286
$ws.Range("K126").Value = 5785.2858
$ws.Range("M126").Value = -3315.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7234.353
$ws.Range("I132").Value = 6198.967
$ws.Range("J132").Value = 14999.75
$ws.Range("K132").Value = 18596.901
$ws.Range("L132").Value = 44999.25
$ws.Range("M132").Value = -16066.901
$ws.Range("N132").Value = -50059.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 49000
$ws.Range("J135").Value = 49000
$ws.Range("L135").Value = 49000
$ws.Range("N135").Value = -59140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5279.7
$ws.Range("I7").Value = 2800.5
$ws.Range("K7").Value = 2800.5
$ws.Range("M7").Value = -2688.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2195.258
$ws.Range("I46").Value = 1437.3846
$ws.Range("J46").Value = 2742.611
$ws.Range("K46").Value = 1437.3846
$ws.Range("L46").Value = 2742.611
$ws.Range("M46").Value = -1249.3846
$ws.Range("N46").Value = -3118.611

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 52298.582
$ws.Range("I100").Value = 79144.39999999999
$ws.Range("J100").Value = 7555.5557
$ws.Range("K100").Value = 79144.39999999999
$ws.Range("L100").Value = 7555.5557
$ws.Range("M100").Value = -78603.39999999999
$ws.Range("N100").Value = -8637.555700000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5298
$ws.Range("I122").Value = 4617.4
$ws.Range("K122").Value = 13852.2
$ws.Range("M122").Value = -11402.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5279.7
$ws.Range("I126").Value = 2800.5
$ws.Range("K126").Value = 8401.5
$ws.Range("M126").Value = -5931.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2847.7058
$ws.Range("I132").Value = 2311.1
$ws.Range("J132").Value = 3614.2856
$ws.Range("K132").Value = 6933.299999999999
$ws.Range("L132").Value = 10842.8568
$ws.Range("M132").Value = -4403.299999999999
$ws.Range("N132").Value = -15902.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1711.8334
$ws.Range("I122").Value = 1675.8125
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5027.4375
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2577.4375
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3748.6
$ws.Range("I132").Value = 921.5
$ws.Range("J132").Value = 5633.3335
$ws.Range("K132").Value = 2764.5
$ws.Range("L132").Value = 16900.0005
$ws.Range("M132").Value = -234.5
$ws.Range("N132").Value = -21960.0005

